# Update the "想去人数" (interested count) figures that changed between
# the two generated-output snapshots.

$wb = $excel.ActiveWorkbook

# Sheet "展览": F2 1334 -> 1335, F3 2852 -> 2855
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 1335
$wsExhibit.Range("F3").Value = 2855

# Sheet "全部类型": F3 1334 -> 1335, F4 2852 -> 2855
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 1335
$wsAll.Range("F4").Value = 2855
